# Append the 26-02-2020 data point as a new row (row 38) to the COVID
# tracking sheet, then move the selection to D38 (mirrors what Excel does
# after typing a new row of data and pressing Tab/Enter off the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(38, 1).Value = "26-02-2020"
$ws.Cells.Item(38, 2).Value = 78191
$ws.Cells.Item(38, 3).Value = 2718

$null = $ws.Range("D38").Select()
